$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2, column B: rename item from "刀" to "Knife"
$ws.Range("B2").Value = "Knife"

# Update the active selection to G2 (matches sheetView selection in the diff)
$ws.Range("G2").Select()
